# Fix #354, #360, #372
#
# plotConfiguration: add "quantiles" / "foldDistance" columns (K/L) and a
# sample "2, 3" foldDistance value for the second plot row.
# plotGrids: add a "tagPrefix" column (D).
# Also updates the active sheet/selection state to match (plotGrids becomes
# the active tab; selections move to the newly added columns).

$wb = $excel.ActiveWorkbook

$wsPlotConfiguration = $wb.Worksheets.Item("plotConfiguration")
$wsPlotConfiguration.Range("K1").Value = "quantiles"
$wsPlotConfiguration.Range("L1").Value = "foldDistance"
$wsPlotConfiguration.Range("L3").Value = "2, 3"
$wsPlotConfiguration.Range("L3").Select() | Out-Null

$wsExportConfiguration = $wb.Worksheets.Item("exportConfiguration")
$wsExportConfiguration.Range("B3").Select() | Out-Null

$wsPlotGrids = $wb.Worksheets.Item("plotGrids")
$wsPlotGrids.Range("D1").Value = "tagPrefix"
$wsPlotGrids.Activate() | Out-Null
$wsPlotGrids.Range("D3").Select() | Out-Null
